$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D (Price) while forcing it to stay a text
# value (some price strings, e.g. "96.82", parse as pure numbers and Excel's
# COM layer would otherwise silently convert them to numeric cells).
function Set-PriceText($row, $text) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-VolumeText($row, $text) {
    $ws.Cells.Item($row, 5).Value = $text
}

# Simple per-row updates: row number, column D (price) new value (or $null if
# unchanged), column E (volume) new value.
$updates = @(
    @{ Row = 2;  D = "43.644.51";  E = "  -0.79%  " },
    @{ Row = 3;  D = "2.293.22";   E = "  +0.73%  " },
    @{ Row = 4;  D = $null;        E = "  -0.02%  " },
    @{ Row = 5;  D = "96.82";      E = "  +1.53%  " },
    @{ Row = 6;  D = "267.18";     E = "  -1.83%  " },
    @{ Row = 7;  D = $null;        E = "  -0.43%  " },
    @{ Row = 8;  D = $null;        E = "  -0.12%  " },
    @{ Row = 9;  D = $null;        E = "  -2.29%  " },
    @{ Row = 10; D = "45.85";      E = "  -2.60%  " },
    @{ Row = 11; D = $null;        E = "  +0.19%  " },
    @{ Row = 12; D = "7.81";       E = "  -3.30%  " },
    @{ Row = 13; D = $null;        E = "  -0.43%  " },
    @{ Row = 14; D = "2.638.87";   E = "  +0.77%  " },
    @{ Row = 15; D = "15.24";      E = "  -1.45%  " },
    @{ Row = 16; D = "0.843";      E = "  +1.95%  " },
    @{ Row = 17; D = "2.287.20";   E = "  +1.19%  " },
    @{ Row = 18; D = "43.598.12";  E = "  -0.99%  " },
    @{ Row = 19; D = $null;        E = "  +2.33%  " },
    @{ Row = 20; D = "6.19";       E = "  +0.07%  " },
    @{ Row = 21; D = "71.91";      E = "  +1.06%  " },
    @{ Row = 22; D = "2.42";       E = "  +5.57%  " },
    @{ Row = 23; D = "232.85";     E = "  -1.59%  " },
    @{ Row = 24; D = "9.20";       E = "  -9.05%  " },
    @{ Row = 25; D = $null;        E = "  +0.10%  " },
    @{ Row = 26; D = "2.50";       E = "  +0.01%  " },
    @{ Row = 27; D = "11.18";      E = "  -1.84%  " },
    @{ Row = 28; D = "3.46";       E = "  +2.95%  " },
    @{ Row = 29; D = "40.36";      E = "  +1.77%  " },
    @{ Row = 30; D = $null;        E = "  +0.52%  " },
    @{ Row = 31; D = "175.57";     E = "  +1.11%  " },
    @{ Row = 32; D = "21.91";      E = "  -1.51%  " },
    @{ Row = 33; D = $null;        E = "  -2.42%  " },
    @{ Row = 34; D = $null;        E = "  -4.33%  " },
    @{ Row = 35; D = $null;        E = "  +1.26%  " },
    @{ Row = 36; D = $null;        E = "  -3.76%  " },
    @{ Row = 37; D = "0.0355";     E = "  +1.07%  " },
    @{ Row = 38; D = "4.38";       E = "  -0.55%  " },
    @{ Row = 39; D = "3.41";       E = "  -1.70%  " },
    @{ Row = 40; D = "0.236";      E = "  -5.16%  " },
    @{ Row = 41; D = $null;        E = "  +4.74%  " },
    @{ Row = 42; D = "12.30";      E = "  -1.00%  " },
    @{ Row = 43; D = "1.34";       E = "  +12.44%  " },
    @{ Row = 44; D = "63.92";      E = "  +3.41%  " },
    @{ Row = 47; D = $null;        E = "  -0.04%  " },
    @{ Row = 48; D = "98.24";      E = "  -1.45%  " },
    @{ Row = 49; D = $null;        E = "  +0.14%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-PriceText $u.Row $u.D
    }
    Set-VolumeText $u.Row $u.E
}

# Rows 45 and 46 swap content (THORChain <-> FraxShare) with updated values.
$ws.Cells.Item(45, 2).Value = "FraxShare"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-PriceText 45 "8.81"
Set-VolumeText 45 "  +2.69%  "

$ws.Cells.Item(46, 2).Value = "THORChain"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-PriceText 46 "5.24"
Set-VolumeText 46 "  -4.38%  "

# Rows 50 and 51 swap content (WOONetwork <-> RocketPoolETH) with updated values.
$ws.Cells.Item(50, 2).Value = "RocketPoolETH"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-PriceText 50 "2.518.48"
Set-VolumeText 50 "  +0.83%  "

$ws.Cells.Item(51, 2).Value = "WOONetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-PriceText 51 "0.429"
Set-VolumeText 51 "  -0.76%  "
